$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card10")

# Remove trailing space from header "Serviced by "
$ws.Range("O1").Value = "Serviced by"

# Fill in the "Serviced by" column for data rows 2-13 with "nan" placeholders,
# matching the rest of the table's empty-value convention.
$ws.Range("O2").Value = "nan"
$ws.Range("O3").Value = "nan"
$ws.Range("O4").Value = "nan"
$ws.Range("O5").Value = "nan"
$ws.Range("O6").Value = "nan"
$ws.Range("O7").Value = "nan"
$ws.Range("O8").Value = "م.محمد عبدالله ،ف.مصطفي"
$ws.Range("O9").Value = "nan"
$ws.Range("O10").Value = "nan"
$ws.Range("O11").Value = "nan"
$ws.Range("O12").Value = "nan"
$ws.Range("O13").Value = "nan"
